# Weekly update: a new week of "Apio" (celery) price data for Terminal La
# Palmera de La Serena was prepended to the historical series. This pushes
# the existing rows 284:306 down by two rows (to 286:308), and the two
# freed-up rows (284 and 285) are populated with the new week's "Primera"
# and "Segunda" quality records.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(284).Insert()
$ws.Rows.Item(284).Insert()

$ws.Cells.Item(284, 1).Value = 8
$ws.Cells.Item(284, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(284, 3).Value = "Coquimbo"
$ws.Cells.Item(284, 4).Value = 44578
$ws.Cells.Item(284, 5).Value = 4
$ws.Cells.Item(284, 6).Value = 100112017
$ws.Cells.Item(284, 7).Value = "Apio"
$ws.Cells.Item(284, 8).Value = "Americana (o)"
$ws.Cells.Item(284, 9).Value = "Primera"
$ws.Cells.Item(284, 10).Value = 2000
$ws.Cells.Item(284, 11).Value = 8000
$ws.Cells.Item(284, 12).Value = 9000
$ws.Cells.Item(284, 13).Value = 8500
$ws.Cells.Item(284, 14).Value = "$/docena de matas"
$ws.Cells.Item(284, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(284, 16).Value = 1417
$ws.Cells.Item(284, 17).Value = 6
$ws.Cells.Item(284, 18).Value = "Hortaliza"

$ws.Cells.Item(285, 1).Value = 8
$ws.Cells.Item(285, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(285, 3).Value = "Coquimbo"
$ws.Cells.Item(285, 4).Value = 44578
$ws.Cells.Item(285, 5).Value = 4
$ws.Cells.Item(285, 6).Value = 100112017
$ws.Cells.Item(285, 7).Value = "Apio"
$ws.Cells.Item(285, 8).Value = "Americana (o)"
$ws.Cells.Item(285, 9).Value = "Segunda"
$ws.Cells.Item(285, 10).Value = 1600
$ws.Cells.Item(285, 11).Value = 6000
$ws.Cells.Item(285, 12).Value = 7000
$ws.Cells.Item(285, 13).Value = 6500
$ws.Cells.Item(285, 14).Value = "$/docena de matas"
$ws.Cells.Item(285, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(285, 16).Value = 1083
$ws.Cells.Item(285, 17).Value = 6
$ws.Cells.Item(285, 18).Value = "Hortaliza"
